# A new daily price record was inserted into the weekly Arándano (blue)
# dataset (Vega Central Mapocho de Santiago), just after the existing row
# for 2022-01-17 (serial 44578). Inserting a whole row shifts every
# following record down by one, which matches the observed diff (every
# row from the old 281 through 332 reappears, unchanged, one row lower,
# and the sheet's used range grows from T332 to T333).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 281..332 down to 282..333, opening up a blank row 281.
$ws.Rows("281:281").Insert()

# Populate the newly opened row with the new record.
$ws.Range("A281").Value = 9
$ws.Range("B281").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C281").Value = 'Metropolitana'
$ws.Range("D281").Value = 45209
$ws.Range("E281").Value = 13
$ws.Range("F281").Value = 'Fruta'
$ws.Range("G281").Value = 100101
$ws.Range("H281").Value = 'Berries'
$ws.Range("I281").Value = 100101001
$ws.Range("J281").Value = 'Arándano (blue)'
$ws.Range("K281").Value = 'Sin especificar'
$ws.Range("L281").Value = 'Primera'
$ws.Range("M281").Value = 160
$ws.Range("N281").Value = 11000
$ws.Range("O281").Value = 12000
$ws.Range("P281").Value = 11500
$ws.Range("Q281").Value = '$/bandeja 2 kilos'
$ws.Range("R281").Value = 'Provincia de Limarí'
$ws.Range("S281").Value = 5750
$ws.Range("T281").Value = 2
